$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("doacoes_registros")
$ws.Activate()

# ---------------------------------------------------------------------------
# 1) Fill in the five template rows (94-98) with the new "2oSPRINT" clothing
#    donation entries. These rows already exist (created earlier as blank
#    "Consistencia semanal" placeholders); we overwrite their contents and
#    adopt the formatting used by the already-filled rows above (90-93).
# ---------------------------------------------------------------------------

# Copy the fully-styled format of row 91 (C:F) onto the C:F cells of the
# rows we are about to populate, so the cell styles match the rest of the
# "2oSPRINT" block instead of the generic empty-row style.
$ws.Range("C91:F91").Copy()
$ws.Range("C94:F98").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$donationRows = @(
    @{ Row = 94; Nome = "Jorge Nazaré";  Grupo = "VIRTUX";              Qtd = 15; Unit = 3 },
    @{ Row = 95; Nome = "Durga";         Grupo = "Motivados Net Supre"; Qtd = 39; Unit = 3 },
    @{ Row = 96; Nome = "Bruno Hudson";  Grupo = "VIRTUX";              Qtd = 30; Unit = 3 },
    @{ Row = 97; Nome = "Bruno Hudson";  Grupo = "VIRTUX";              Qtd = 30; Unit = 3 },
    @{ Row = 98; Nome = "Bruno Hudson";  Grupo = "VIRTUX";              Qtd = 12; Unit = 3 }
)

foreach ($entry in $donationRows) {
    $r = $entry.Row
    $ws.Range("A$r").Value = "2ºSPRINT"
    $ws.Range("B$r").Value = 45987
    $ws.Range("C$r").Value = $entry.Nome
    $ws.Range("D$r").Value = $entry.Grupo
    $ws.Range("E$r").Value = "Roupas"
    $ws.Range("F$r").Value = "Peça de roupa"
    $ws.Range("G$r").Value = $entry.Qtd
    $ws.Range("H$r").Value = $entry.Unit
    $ws.Range("L$r").Value = "Roupas"
}

# Row 99 stays an (otherwise empty) template row, but its SPRINT tag moves
# from "Consistencia semanal" to "2ºSPRINT" along with the block above.
$ws.Range("A99").Value = "2ºSPRINT"

# ---------------------------------------------------------------------------
# 2) Insert six fresh blank template rows before the "3oSPRINT"/"ENCERRAMENTO"
#    rows (old rows 107-111 shift down to 113-117).
# ---------------------------------------------------------------------------
$ws.Rows("107:112").Insert(-4121, 0)

# The inserted rows come back with a generic/default style; restore the
# borders + alignment used by the rest of the blank template rows
# (s=9 general cols, s=10 date col, s=11 name col, s=12 group/item cols,
#  s=17 category col) so the new rows match rows 100-106 exactly.
for ($r = 107; $r -le 112; $r++) {
    $general = $ws.Range("A$r,G$r,H$r,I$r,J$r,K$r,L$r")
    $general.Borders.LineStyle = 1
    $general.HorizontalAlignment = -4108

    $dateCell = $ws.Range("B$r")
    $dateCell.Borders.LineStyle = 1
    $dateCell.HorizontalAlignment = -4108
    $dateCell.NumberFormat = "dd/mm/yy"

    $nameCell = $ws.Range("C$r")
    $nameCell.Borders.LineStyle = 1

    $groupItemCells = $ws.Range("D$r,F$r")
    $groupItemCells.Borders.LineStyle = 1
    $groupItemCells.HorizontalAlignment = -4131

    $categoryCell = $ws.Range("E$r")
    $categoryCell.Borders.LineStyle = 1
    $categoryCell.Borders.Item(7).LineStyle = -4142
    $categoryCell.Borders.Item(10).LineStyle = -4142
    $categoryCell.HorizontalAlignment = -4131

    $ws.Range("I$r").Formula = "=G$r*H$r"
    $ws.Range("K$r").Formula = "=I$r+J$r"
}

# ---------------------------------------------------------------------------
# 3) Move the active selection to reflect where the author ended up editing.
# ---------------------------------------------------------------------------
$ws.Range("G95").Select()
